# Generate Report for Handback
# Simulates the "handback" CI step: the de-de locale has just been handed
# back (localized xliff processed) and both zh-cn/de-de rows pick up their
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# values, the overview status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and a couple of report columns are
# widened so the new long file names are readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16289712e6e111d66615baddb73fdbc2b06a989f/e2e/"
$aMdUrl = $githubBase + "a.md"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn / de-de: fill in Latest Target File (I), Latest Handback File
#    (J) and Latest Handback DateTime (K) for both data rows.
# ---------------------------------------------------------------------

# --- zh-cn ---
$wsZhCn.Range("I2").Value = "a.md"
$wsZhCn.Range("I2").Style = "Hyperlink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $aMdUrl, "", "", "a.md")

$wsZhCn.Range("I3").Value = "a.md"
$wsZhCn.Range("I3").Style = "Hyperlink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $aMdUrl, "", "", "a.md")

$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-09-04 18:41:38"
$wsZhCn.Range("K3").Value = "2016-09-04 18:41:38"

# --- de-de ---
$wsDeDe.Range("I2").Value = "a.md"
$wsDeDe.Range("I2").Style = "Hyperlink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $aMdUrl, "", "", "a.md")

$wsDeDe.Range("I3").Value = "a.md"
$wsDeDe.Range("I3").Style = "Hyperlink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $aMdUrl, "", "", "a.md")

$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-09-04 18:41:46"
$wsDeDe.Range("K3").Value = "2016-09-04 18:41:46"

# ---------------------------------------------------------------------
# 3. Widen columns so the new, longer values are readable.
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.1   # zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = 29.1   # de-de status column

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1    # Status
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15  # Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1    # Status
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15  # Latest Handback File
